# Update the title-slide placeholders on slide 1:
#   "Insert Title" -> "Loan Analysis and Prediction"
#   "Insert Date"  -> "2022"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleShape = $s.Shapes.Item("Text Placeholder 2")
$titleShape.TextFrame.TextRange.Text = "Loan Analysis and Prediction"

$dateShape = $s.Shapes.Item("Text Placeholder 4")
$dateShape.TextFrame.TextRange.Text = "2022"
